# Auto-generated Excel COM-interop script to apply the Typhon_Profits data refresh
# Updates currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ
# columns (H-N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 76: H76=3971425.2, I76=3383.125, K76=3383.125, M76=-3068.125
$ws.Range("H76").Value = 3971425.2
$ws.Range("I76").Value = 3383.125
$ws.Range("K76").Value = 3383.125
$ws.Range("M76").Value = -3068.125
# Row 79: H79=3971425.2, I79=3383.125, K79=3383.125, M79=-2291.125
$ws.Range("H79").Value = 3971425.2
$ws.Range("I79").Value = 3383.125
$ws.Range("K79").Value = 3383.125
$ws.Range("M79").Value = -2291.125
# Row 98: H98=1153.4706, I98=510, K98=510, M98=988
$ws.Range("H98").Value = 1153.4706
$ws.Range("I98").Value = 510
$ws.Range("K98").Value = 510
$ws.Range("M98").Value = 988
# Row 122: H122=1153.4706, I122=510, K122=1530, M122=920
$ws.Range("H122").Value = 1153.4706
$ws.Range("I122").Value = 510
$ws.Range("K122").Value = 1530
$ws.Range("M122").Value = 920
# Row 135: H135=62501904, I135=2240, J135=166668020, K135=20160, L135=1500012180, M135=-17625, N135=-1500017250
$ws.Range("H135").Value = 62501904
$ws.Range("I135").Value = 2240
$ws.Range("J135").Value = 166668020
$ws.Range("K135").Value = 20160
$ws.Range("L135").Value = 1500012180
$ws.Range("M135").Value = -17625
$ws.Range("N135").Value = -1500017250
# Row 137: H137=1972.5834, I137=1740.8889, K137=5222.6667, M137=-2672.6667
$ws.Range("H137").Value = 1972.5834
$ws.Range("I137").Value = 1740.8889
$ws.Range("K137").Value = 5222.6667
$ws.Range("M137").Value = -2672.6667
# Row 138: H138=34486100, J138=4589.7856, L138=13769.3568, N138=-24049.3568
$ws.Range("H138").Value = 34486100
$ws.Range("J138").Value = 4589.7856
$ws.Range("L138").Value = 13769.3568
$ws.Range("N138").Value = -24049.3568
# Row 141: H141=2291, I141=1538, K141=4614, M141=566
$ws.Range("H141").Value = 2291
$ws.Range("I141").Value = 1538
$ws.Range("K141").Value = 4614
$ws.Range("M141").Value = 566

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32: H32=3892.4285, I32=3037.8, J32=15002.6, K32=3037.8, L32=15002.6, M32=-2750.8, N32=-15576.6
$ws.Range("H32").Value = 3892.4285
$ws.Range("I32").Value = 3037.8
$ws.Range("J32").Value = 15002.6
$ws.Range("K32").Value = 3037.8
$ws.Range("L32").Value = 15002.6
$ws.Range("M32").Value = -2750.8
$ws.Range("N32").Value = -15576.6
# Row 61: H61=1980.9688, I61=1496.1072, K61=1496.1072, M61=-1284.1072
$ws.Range("H61").Value = 1980.9688
$ws.Range("I61").Value = 1496.1072
$ws.Range("K61").Value = 1496.1072
$ws.Range("M61").Value = -1284.1072
# Row 74: H74=55557656, I74=71429630, J74=5725, K74=71429630, L74=5725, M74=-71428756, N74=-7473
$ws.Range("H74").Value = 55557656
$ws.Range("I74").Value = 71429630
$ws.Range("J74").Value = 5725
$ws.Range("K74").Value = 71429630
$ws.Range("L74").Value = 5725
$ws.Range("M74").Value = -71428756
$ws.Range("N74").Value = -7473
# Row 77: H77=55557656, I77=71429630, J77=5725, K77=357148150, L77=28625, M77=-357143782, N77=-37361
$ws.Range("H77").Value = 55557656
$ws.Range("I77").Value = 71429630
$ws.Range("J77").Value = 5725
$ws.Range("K77").Value = 357148150
$ws.Range("L77").Value = 28625
$ws.Range("M77").Value = -357143782
$ws.Range("N77").Value = -37361
# Row 97: H97=90910260, I97=1472.5, K97=1472.5, M97=-976.5
$ws.Range("H97").Value = 90910260
$ws.Range("I97").Value = 1472.5
$ws.Range("K97").Value = 1472.5
$ws.Range("M97").Value = -976.5
# Row 132: H132=14045.075, I132=1201.7931, J132=47904.637, K132=3605.379300000001, L132=143713.911, M132=-1075.379300000001, N132=-148773.911
$ws.Range("H132").Value = 14045.075
$ws.Range("I132").Value = 1201.7931
$ws.Range("J132").Value = 47904.637
$ws.Range("K132").Value = 3605.379300000001
$ws.Range("L132").Value = 143713.911
$ws.Range("M132").Value = -1075.379300000001
$ws.Range("N132").Value = -148773.911
# Row 136: H136=1980.9688, I136=1496.1072, K136=4488.321599999999, M136=-1938.321599999999
$ws.Range("H136").Value = 1980.9688
$ws.Range("I136").Value = 1496.1072
$ws.Range("K136").Value = 4488.321599999999
$ws.Range("M136").Value = -1938.321599999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 105: H105=3574968.5, I105=4182.5, J105=8336016.5, K105=4182.5, L105=8336016.5, M105=-2435.5, N105=-8339510.5
$ws.Range("H105").Value = 3574968.5
$ws.Range("I105").Value = 4182.5
$ws.Range("J105").Value = 8336016.5
$ws.Range("K105").Value = 4182.5
$ws.Range("L105").Value = 8336016.5
$ws.Range("M105").Value = -2435.5
$ws.Range("N105").Value = -8339510.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31: H31=3043.75, I31=2502.4375, K31=2502.4375, M31=-2207.4375
$ws.Range("H31").Value = 3043.75
$ws.Range("I31").Value = 2502.4375
$ws.Range("K31").Value = 2502.4375
$ws.Range("M31").Value = -2207.4375
# Row 34: H34=3043.75, I34=2502.4375, K34=2502.4375, M34=-2300.4375
$ws.Range("H34").Value = 3043.75
$ws.Range("I34").Value = 2502.4375
$ws.Range("K34").Value = 2502.4375
$ws.Range("M34").Value = -2300.4375
# Row 58: H58=19571.215, I58=1684.6666, K58=1684.6666, M58=-1481.6666
$ws.Range("H58").Value = 19571.215
$ws.Range("I58").Value = 1684.6666
$ws.Range("K58").Value = 1684.6666
$ws.Range("M58").Value = -1481.6666
# Row 132: H132=3333.682, I132=2664.2, K132=7992.599999999999, M132=-5462.599999999999
$ws.Range("H132").Value = 3333.682
$ws.Range("I132").Value = 2664.2
$ws.Range("K132").Value = 7992.599999999999
$ws.Range("M132").Value = -5462.599999999999
# Row 136: H136=19571.215, I136=1684.6666, K136=5053.9998, M136=-2503.9998
$ws.Range("H136").Value = 19571.215
$ws.Range("I136").Value = 1684.6666
$ws.Range("K136").Value = 5053.9998
$ws.Range("M136").Value = -2503.9998

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 87: H87=11066.533, I87=4590.727, J87=28875, K87=13772.181, L87=86625, M87=-12524.181, N87=-89121
$ws.Range("H87").Value = 11066.533
$ws.Range("I87").Value = 4590.727
$ws.Range("J87").Value = 28875
$ws.Range("K87").Value = 13772.181
$ws.Range("L87").Value = 86625
$ws.Range("M87").Value = -12524.181
$ws.Range("N87").Value = -89121
# Row 90: H90=11066.533, I90=4590.727, J90=28875, K90=41316.543, L90=259875, M90=-35076.543, N90=-272355
$ws.Range("H90").Value = 11066.533
$ws.Range("I90").Value = 4590.727
$ws.Range("J90").Value = 28875
$ws.Range("K90").Value = 41316.543
$ws.Range("L90").Value = 259875
$ws.Range("M90").Value = -35076.543
$ws.Range("N90").Value = -272355
# Row 92: H92=428, I92=260.2, J92=847.5, K92=780.5999999999999, L92=2542.5, M92=467.4000000000001, N92=-5038.5
$ws.Range("H92").Value = 428
$ws.Range("I92").Value = 260.2
$ws.Range("J92").Value = 847.5
$ws.Range("K92").Value = 780.5999999999999
$ws.Range("L92").Value = 2542.5
$ws.Range("M92").Value = 467.4000000000001
$ws.Range("N92").Value = -5038.5
# Row 114: H114=806.5, J114=1115, L114=3345, N114=-9853
$ws.Range("H114").Value = 806.5
$ws.Range("J114").Value = 1115
$ws.Range("L114").Value = 3345
$ws.Range("N114").Value = -9853

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 28: H28=0, J28=0, L28=0
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
# Row 70: H70=12522520, J70=12522520, L70=12522520, N70=-12523060
$ws.Range("H70").Value = 12522520
$ws.Range("J70").Value = 12522520
$ws.Range("L70").Value = 12522520
$ws.Range("N70").Value = -12523060
# Row 73: H73=12522520, J73=12522520, L73=12522520, N73=-12524392
$ws.Range("H73").Value = 12522520
$ws.Range("J73").Value = 12522520
$ws.Range("L73").Value = 12522520
$ws.Range("N73").Value = -12524392
# Row 97: H97=1323, I97=1461.4667, J97=976.8333, K97=1461.4667, L97=976.8333, M97=-965.4666999999999, N97=-1968.8333
$ws.Range("H97").Value = 1323
$ws.Range("I97").Value = 1461.4667
$ws.Range("J97").Value = 976.8333
$ws.Range("K97").Value = 1461.4667
$ws.Range("L97").Value = 976.8333
$ws.Range("M97").Value = -965.4666999999999
$ws.Range("N97").Value = -1968.8333
# Row 113: H113=2208.3635, I113=1739.8572, K113=1739.8572, M113=430.1428000000001
$ws.Range("H113").Value = 2208.3635
$ws.Range("I113").Value = 1739.8572
$ws.Range("K113").Value = 1739.8572
$ws.Range("M113").Value = 430.1428000000001
# Row 132: H132=23802.783, I132=1622.2858, J132=58305.777, K132=4866.857400000001, L132=174917.331, M132=-2336.857400000001, N132=-179977.331
$ws.Range("H132").Value = 23802.783
$ws.Range("I132").Value = 1622.2858
$ws.Range("J132").Value = 58305.777
$ws.Range("K132").Value = 4866.857400000001
$ws.Range("L132").Value = 174917.331
$ws.Range("M132").Value = -2336.857400000001
$ws.Range("N132").Value = -179977.331
# Row 28: clear N28 entirely (no longer present after refresh)
$ws.Range("N28").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 42: H42=4000, J42=4000, L42=4000, N42=-5126
$ws.Range("H42").Value = 4000
$ws.Range("J42").Value = 4000
$ws.Range("L42").Value = 4000
$ws.Range("N42").Value = -5126
# Row 49: H49=4000, J49=4000, L49=4000, N49=-4294
$ws.Range("H49").Value = 4000
$ws.Range("J49").Value = 4000
$ws.Range("L49").Value = 4000
$ws.Range("N49").Value = -4294
# Row 68: H68=2297.3157, I68=2260.1, K68=2260.1, M68=-1511.1
$ws.Range("H68").Value = 2297.3157
$ws.Range("I68").Value = 2260.1
$ws.Range("K68").Value = 2260.1
$ws.Range("M68").Value = -1511.1
# Row 71: H71=2297.3157, I71=2260.1, K71=11300.5, M71=-7556.5
$ws.Range("H71").Value = 2297.3157
$ws.Range("I71").Value = 2260.1
$ws.Range("K71").Value = 11300.5
$ws.Range("M71").Value = -7556.5
# Row 93: H93=1111, I93=846.7059, J93=1859.8334, K93=846.7059, L93=1859.8334, M93=401.2941, N93=-4355.8334
$ws.Range("H93").Value = 1111
$ws.Range("I93").Value = 846.7059
$ws.Range("J93").Value = 1859.8334
$ws.Range("K93").Value = 846.7059
$ws.Range("L93").Value = 1859.8334
$ws.Range("M93").Value = 401.2941
$ws.Range("N93").Value = -4355.8334
# Row 132: H132=711160.5600000001, I132=1508029, J132=2833.111, K132=4524087, L132=8499.332999999999, M132=-4521557, N132=-13559.333
$ws.Range("H132").Value = 711160.5600000001
$ws.Range("I132").Value = 1508029
$ws.Range("J132").Value = 2833.111
$ws.Range("K132").Value = 4524087
$ws.Range("L132").Value = 8499.332999999999
$ws.Range("M132").Value = -4521557
$ws.Range("N132").Value = -13559.333
# Row 136: H136=0, I136=0, K136=0
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
# Row 136: clear M136 entirely (no longer present after refresh)
$ws.Range("M136").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132: H132=1649.1538, I132=1239.8334, J132=2000, K132=3719.5002, L132=6000, M132=-1189.5002, N132=-11060
$ws.Range("H132").Value = 1649.1538
$ws.Range("I132").Value = 1239.8334
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 3719.5002
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -1189.5002
$ws.Range("N132").Value = -11060
# Row 136: H136=33300646, I136=44882390, J136=3125.625, K136=134647170, L136=9376.875, M136=-134644620, N136=-14476.875
$ws.Range("H136").Value = 33300646
$ws.Range("I136").Value = 44882390
$ws.Range("J136").Value = 3125.625
$ws.Range("K136").Value = 134647170
$ws.Range("L136").Value = 9376.875
$ws.Range("M136").Value = -134644620
$ws.Range("N136").Value = -14476.875
